$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) to lowercase/snake_case column names
$ws.Range("A1").Value = "port_name"
$ws.Range("B1").Value = "material"
$ws.Range("C1").Value = "handling_cost_inr_tonne"
$ws.Range("D1").Value = "storage_cost_inr_tonne_day"
$ws.Range("E1").Value = "max_throughput_t_day"

# Update data values (columns C, D, E) for rows 2-13
$ws.Cells.Item(2, 3).Value = 213.48
$ws.Cells.Item(2, 4).Value = 11.6
$ws.Cells.Item(2, 5).Value = 26900

$ws.Cells.Item(3, 3).Value = 197.44
$ws.Cells.Item(3, 4).Value = 7.45
$ws.Cells.Item(3, 5).Value = 20300

$ws.Cells.Item(4, 3).Value = 303.12
$ws.Cells.Item(4, 4).Value = 13.06
$ws.Cells.Item(4, 5).Value = 11900

$ws.Cells.Item(5, 3).Value = 318.08
$ws.Cells.Item(5, 4).Value = 9.57
$ws.Cells.Item(5, 5).Value = 22800

$ws.Cells.Item(6, 3).Value = 260.67
$ws.Cells.Item(6, 4).Value = 5.95
$ws.Cells.Item(6, 5).Value = 17400

$ws.Cells.Item(7, 3).Value = 291.21
$ws.Cells.Item(7, 4).Value = 11.83
$ws.Cells.Item(7, 5).Value = 11800

$ws.Cells.Item(8, 3).Value = 332.47
$ws.Cells.Item(8, 4).Value = 9.619999999999999
$ws.Cells.Item(8, 5).Value = 21100

$ws.Cells.Item(9, 3).Value = 222.09
$ws.Cells.Item(9, 4).Value = 6.03
$ws.Cells.Item(9, 5).Value = 18100

$ws.Cells.Item(10, 3).Value = 419.1
$ws.Cells.Item(10, 4).Value = 11.03
$ws.Cells.Item(10, 5).Value = 13900

$ws.Cells.Item(11, 3).Value = 235.36
$ws.Cells.Item(11, 4).Value = 10.51
$ws.Cells.Item(11, 5).Value = 21900

$ws.Cells.Item(12, 3).Value = 210.42
$ws.Cells.Item(12, 4).Value = 5.35
$ws.Cells.Item(12, 5).Value = 16600

$ws.Cells.Item(13, 3).Value = 376.55
$ws.Cells.Item(13, 4).Value = 14.96
$ws.Cells.Item(13, 5).Value = 17000

$wb.Save()
